# "Minor corrections and adds"
#
# The source data correction: row 10 on "Приходы" (Income) represented an
# employee with 38 sick/absence days (column D) and an incorrectly entered
# "E" adjustment of 2 days that shouldn't have been there. Clearing E10
# removes that erroneous adjustment so G10 (=C10*(D10-E10)) correctly
# computes off the full 38 instead of 36.
#
# Every other number touched by the original diff (G10, G20/G21 totals on
# "Приходы"; C10/E10/F10 on "Итоги"; the chart caches fed by those ranges)
# is a formula result that recalculates automatically off this single
# input edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Приходы")

# Remove the stray adjustment value that was in E10.
$ws.Range("E10").ClearContents()
